{"js": "// Remove the trailing \"Ver no Jupiter ...\" line, the \"\u00a9 2020 ...\" footer\n// line that follows it, and the blank paragraph that used to separate that\n// footer from the final (page-break) paragraph. The blank paragraph that\n// sits directly after the bibliography entry is left untouched.\nconst body = context.document.body;\n\n// Locate the \"Ver no Jupiter...\" paragraph by its text so the edit is not\n// dependent on a hard-coded paragraph index.\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const jupiterPara = results.items[0].paragraphs.getFirst();\n  const copyrightPara = jupiterPara.getNext();\n  const trailingBlankPara = copyrightPara.getNext();\n\n  // Delete from the end backwards so earlier handles stay valid.\n  trailingBlankPara.delete();\n  copyrightPara.delete();\n  jupiterPara.delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" line, the \"\u00a9 2020 ...\" footer\n# line that follows it, and the blank paragraph that used to separate that\n# footer from the final (page-break) paragraph. The blank paragraph that\n# sits directly after the bibliography entry is left untouched.\n$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter...\" paragraph using Find rather than a\n# hard-coded paragraph index.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\n\nif ($found) {\n    $all = $d.Paragraphs\n    $targetIndex = $null\n    for ($i = 1; $i -le $all.Count; $i++) {\n        $p = $all.Item($i)\n        if ($p.Range.Start -le $findRange.Start -and $p.Range.End -ge $findRange.End) {\n            $targetIndex = $i\n            break\n        }\n    }\n\n    if ($targetIndex -ne $null) {\n        $jupiterPara = $d.Paragraphs.Item($targetIndex)\n        $copyrightPara = $jupiterPara.Next()\n        $trailingBlankPara = $copyrightPara.Next()\n\n        $deleteRange = $d.Range($jupiterPara.Range.Start, $trailingBlankPara.Range.End)\n        $deleteRange.Delete()\n    }\n}\n"}
